$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1433.3334
$ws.Range("I19").Value = 1121.2
$ws.Range("J19").Value = 2057.6
$ws.Range("K19").Value = 1121.2
$ws.Range("L19").Value = 2057.6
$ws.Range("M19").Value = -946.2
$ws.Range("N19").Value = -2407.6
$ws.Range("H43").Value = 1026275
$ws.Range("I43").Value = 2100
$ws.Range("J43").Value = 1367666.6
$ws.Range("K43").Value = 2100
$ws.Range("L43").Value = 1367666.6
$ws.Range("M43").Value = -2031
$ws.Range("N43").Value = -1367804.6
$ws.Range("H49").Value = 334.25
$ws.Range("I49").Value = 550
$ws.Range("J49").Value = 310.27777
$ws.Range("K49").Value = 1650
$ws.Range("L49").Value = 930.83331
$ws.Range("M49").Value = -1514
$ws.Range("N49").Value = -1202.83331
$ws.Range("H58").Value = 35717532
$ws.Range("J58").Value = 100007560
$ws.Range("L58").Value = 300022680
$ws.Range("N58").Value = -300022980
$ws.Range("H80").Value = 26753.264
$ws.Range("I80").Value = 7937.385
$ws.Range("J80").Value = 67521
$ws.Range("K80").Value = 23812.155
$ws.Range("L80").Value = 202563
$ws.Range("M80").Value = -22814.155
$ws.Range("N80").Value = -204559
$ws.Range("H83").Value = 26753.264
$ws.Range("I83").Value = 7937.385
$ws.Range("J83").Value = 67521
$ws.Range("K83").Value = 71436.465
$ws.Range("L83").Value = 607689
$ws.Range("M83").Value = -66444.465
$ws.Range("N83").Value = -617673
$ws.Range("H103").Value = 493.64285
$ws.Range("I103").Value = 622.2
$ws.Range("J103").Value = 453.46875
$ws.Range("K103").Value = 1866.6
$ws.Range("L103").Value = 1360.40625
$ws.Range("M103").Value = -1280.6
$ws.Range("N103").Value = -2532.40625
$ws.Range("H116").Value = 20845748
$ws.Range("I116").Value = 62505990
$ws.Range("J116").Value = 15625.75
$ws.Range("K116").Value = 62505990
$ws.Range("L116").Value = 15625.75
$ws.Range("M116").Value = -62502548
$ws.Range("N116").Value = -22509.75
$ws.Range("H138").Value = 3509.8298
$ws.Range("J138").Value = 3947.0144
$ws.Range("L138").Value = 11841.0432
$ws.Range("N138").Value = -22121.0432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 41752.08
$ws.Range("I74").Value = 54873.895
$ws.Range("J74").Value = 6135.7144
$ws.Range("K74").Value = 54873.895
$ws.Range("L74").Value = 6135.7144
$ws.Range("M74").Value = -53999.895
$ws.Range("N74").Value = -7883.7144
$ws.Range("H77").Value = 41752.08
$ws.Range("I77").Value = 54873.895
$ws.Range("J77").Value = 6135.7144
$ws.Range("K77").Value = 274369.475
$ws.Range("L77").Value = 30678.572
$ws.Range("M77").Value = -270001.475
$ws.Range("N77").Value = -39414.572
$ws.Range("H97").Value = 5953907
$ws.Range("I97").Value = 959.1111
$ws.Range("K97").Value = 959.1111
$ws.Range("M97").Value = -463.1111
$ws.Range("H101").Value = 40956.332
$ws.Range("J101").Value = 40956.332
$ws.Range("L101").Value = 40956.332
$ws.Range("N101").Value = -47446.332
$ws.Range("H102").Value = 12504650
$ws.Range("I102").Value = 13337626
$ws.Range("K102").Value = 13337626
$ws.Range("M102").Value = -13336004

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 22200
$ws.Range("J21").Value = 22200
$ws.Range("L21").Value = 22200
$ws.Range("N21").Value = -22672
$ws.Range("H26").Value = 11825.857
$ws.Range("I26").Value = 4823.1665
$ws.Range("J26").Value = 53842
$ws.Range("K26").Value = 4823.1665
$ws.Range("L26").Value = 53842
$ws.Range("M26").Value = -4531.1665
$ws.Range("N26").Value = -54426
$ws.Range("H27").Value = 53642
$ws.Range("J27").Value = 53642
$ws.Range("L27").Value = 53642
$ws.Range("N27").Value = -54026
$ws.Range("H28").Value = 41921
$ws.Range("J28").Value = 41921
$ws.Range("L28").Value = 41921
$ws.Range("N28").Value = -42509
$ws.Range("H33").Value = 7283.625
$ws.Range("J33").Value = 8000
$ws.Range("L33").Value = 8000
$ws.Range("N33").Value = -8672

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10203.808
$ws.Range("I31").Value = 3480.111
$ws.Range("J31").Value = 13763.412
$ws.Range("K31").Value = 3480.111
$ws.Range("L31").Value = 13763.412
$ws.Range("M31").Value = -3185.111
$ws.Range("N31").Value = -14353.412
$ws.Range("H34").Value = 10203.808
$ws.Range("I34").Value = 3480.111
$ws.Range("J34").Value = 13763.412
$ws.Range("K34").Value = 3480.111
$ws.Range("L34").Value = 13763.412
$ws.Range("M34").Value = -3278.111
$ws.Range("N34").Value = -14167.412
$ws.Range("H58").Value = 5267.0176
$ws.Range("I58").Value = 2460.8157
$ws.Range("K58").Value = 2460.8157
$ws.Range("M58").Value = -2257.8157
$ws.Range("H103").Value = 49921
$ws.Range("J103").Value = 53842
$ws.Range("L103").Value = 53842
$ws.Range("N103").Value = -56186
$ws.Range("H132").Value = 9443.611000000001
$ws.Range("I132").Value = 6141.4287
$ws.Range("J132").Value = 11545
$ws.Range("K132").Value = 18424.2861
$ws.Range("L132").Value = 34635
$ws.Range("M132").Value = -15894.2861
$ws.Range("N132").Value = -39695
$ws.Range("H134").Value = 6264.636
$ws.Range("J134").Value = 10885.4
$ws.Range("L134").Value = 32656.2
$ws.Range("N134").Value = -37726.2
$ws.Range("H136").Value = 5267.0176
$ws.Range("I136").Value = 2460.8157
$ws.Range("K136").Value = 7382.4471
$ws.Range("M136").Value = -4832.4471

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2015.8334
$ws.Range("I5").Value = 903.3333
$ws.Range("K5").Value = 2709.9999
$ws.Range("M5").Value = -2597.9999
$ws.Range("H11").Value = 2399.6
$ws.Range("I11").Value = 49
$ws.Range("J11").Value = 3966.6667
$ws.Range("K11").Value = 147
$ws.Range("L11").Value = 11900.0001
$ws.Range("M11").Value = -7
$ws.Range("N11").Value = -12180.0001
$ws.Range("H12").Value = 4130.857
$ws.Range("J12").Value = 4113.3335
$ws.Range("L12").Value = 12340.0005
$ws.Range("N12").Value = -12686.0005
$ws.Range("H107").Value = 4681551.5
$ws.Range("I107").Value = 1428987.9
$ws.Range("J107").Value = 5250750.5
$ws.Range("K107").Value = 4286963.699999999
$ws.Range("L107").Value = 15752251.5
$ws.Range("M107").Value = -4285043.699999999
$ws.Range("N107").Value = -15756091.5
$ws.Range("H114").Value = 1975
$ws.Range("I114").Value = 462.5
$ws.Range("J114").Value = 5000
$ws.Range("K114").Value = 1387.5
$ws.Range("L114").Value = 15000
$ws.Range("M114").Value = 1866.5
$ws.Range("N114").Value = -21508
$ws.Range("H132").Value = 8789.162
$ws.Range("I132").Value = 4984.316
$ws.Range("J132").Value = 12805.389
$ws.Range("K132").Value = 44858.844
$ws.Range("L132").Value = 115248.501
$ws.Range("M132").Value = -42328.844
$ws.Range("N132").Value = -120308.501
$ws.Range("H135").Value = 2015.8334
$ws.Range("I135").Value = 903.3333
$ws.Range("K135").Value = 8129.9997
$ws.Range("M135").Value = -5594.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4469.5
$ws.Range("I102").Value = 4342.4165
$ws.Range("K102").Value = 4342.4165
$ws.Range("M102").Value = -2720.4165
$ws.Range("H122").Value = 4026155
$ws.Range("I122").Value = 4262870
$ws.Range("J122").Value = 1998
$ws.Range("K122").Value = 12788610
$ws.Range("L122").Value = 5994
$ws.Range("M122").Value = -12786160
$ws.Range("N122").Value = -10894
$ws.Range("H126").Value = 8356.081
$ws.Range("I126").Value = 4797.273
$ws.Range("J126").Value = 9861.73
$ws.Range("K126").Value = 14391.819
$ws.Range("L126").Value = 29585.19
$ws.Range("M126").Value = -11921.819
$ws.Range("N126").Value = -34525.19
$ws.Range("H132").Value = 5224.4688
$ws.Range("I132").Value = 2976.4736
$ws.Range("J132").Value = 8510
$ws.Range("K132").Value = 8929.4208
$ws.Range("L132").Value = 25530
$ws.Range("M132").Value = -6399.4208
$ws.Range("N132").Value = -30590

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7855.4443
$ws.Range("I40").Value = 5566.3335
$ws.Range("K40").Value = 5566.3335
$ws.Range("M40").Value = -5430.3335
$ws.Range("H82").Value = 39346.777
$ws.Range("I82").Value = 79229.16
$ws.Range("K82").Value = 79229.16
$ws.Range("M82").Value = -78868.16
$ws.Range("H85").Value = 39346.777
$ws.Range("I85").Value = 79229.16
$ws.Range("K85").Value = 79229.16
$ws.Range("M85").Value = -77981.16
$ws.Range("H136").Value = 14353.842
$ws.Range("I136").Value = 4188.7144
$ws.Range("J136").Value = 20283.5
$ws.Range("K136").Value = 12566.1432
$ws.Range("L136").Value = 60850.5
$ws.Range("M136").Value = -10016.1432
$ws.Range("N136").Value = -65950.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9677.4
$ws.Range("I132").Value = 14311.556
$ws.Range("J132").Value = 5885.8184
$ws.Range("K132").Value = 42934.66800000001
$ws.Range("L132").Value = 17657.4552
$ws.Range("M132").Value = -40404.66800000001
$ws.Range("N132").Value = -22717.4552
